# Add "position/" prefix to every "position_<n>.png" image filename
# stored in column A of Sheet1 (rows 2 through the last used row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "position_*.png" -and $val -notlike "position/*") {
        $cell.Value2 = "position/" + $val
    }
}
